{"js": "// Update the two-digit-by-two-digit multiplication prompts in the practice\n// table. Each cell's text (e.g. \"45\u00d736=\") is unique in the document, so we\n// can safely locate and replace each one independently via search().\nconst replacements = [\n  [\"45\u00d736=\", \"76\u00d786=\"],\n  [\"19\u00d756=\", \"95\u00d766=\"],\n  [\"88\u00d765=\", \"53\u00d769=\"],\n  [\"74\u00d783=\", \"67\u00d790=\"],\n  [\"49\u00d799=\", \"35\u00d770=\"],\n  [\"95\u00d786=\", \"54\u00d767=\"],\n  [\"23\u00d781=\", \"53\u00d785=\"],\n  [\"49\u00d729=\", \"95\u00d729=\"],\n  [\"37\u00d711=\", \"26\u00d759=\"],\n  [\"94\u00d756=\", \"79\u00d790=\"],\n  [\"99\u00d777=\", \"38\u00d738=\"],\n  [\"89\u00d741=\", \"36\u00d783=\"],\n  [\"19\u00d762=\", \"20\u00d747=\"],\n  [\"42\u00d717=\", \"14\u00d797=\"],\n  [\"88\u00d784=\", \"86\u00d762=\"],\n  [\"55\u00d774=\", \"82\u00d723=\"],\n  [\"12\u00d793=\", \"97\u00d768=\"],\n  [\"72\u00d760=\", \"90\u00d797=\"],\n  [\"27\u00d728=\", \"83\u00d781=\"],\n  [\"61\u00d755=\", \"59\u00d750=\"],\n  [\"74\u00d736=\", \"66\u00d786=\"],\n  [\"36\u00d797=\", \"97\u00d748=\"],\n  [\"30\u00d711=\", \"78\u00d744=\"],\n  [\"44\u00d720=\", \"93\u00d717=\"],\n  [\"72\u00d755=\", \"86\u00d712=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit-by-two-digit multiplication prompts in the practice\n# table. Each cell's text (e.g. \"45\u00d736=\") is unique in the document, so each\n# Find/Replace pair targets exactly one cell.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"45\u00d736=\"; New = \"76\u00d786=\" },\n    @{ Old = \"19\u00d756=\"; New = \"95\u00d766=\" },\n    @{ Old = \"88\u00d765=\"; New = \"53\u00d769=\" },\n    @{ Old = \"74\u00d783=\"; New = \"67\u00d790=\" },\n    @{ Old = \"49\u00d799=\"; New = \"35\u00d770=\" },\n    @{ Old = \"95\u00d786=\"; New = \"54\u00d767=\" },\n    @{ Old = \"23\u00d781=\"; New = \"53\u00d785=\" },\n    @{ Old = \"49\u00d729=\"; New = \"95\u00d729=\" },\n    @{ Old = \"37\u00d711=\"; New = \"26\u00d759=\" },\n    @{ Old = \"94\u00d756=\"; New = \"79\u00d790=\" },\n    @{ Old = \"99\u00d777=\"; New = \"38\u00d738=\" },\n    @{ Old = \"89\u00d741=\"; New = \"36\u00d783=\" },\n    @{ Old = \"19\u00d762=\"; New = \"20\u00d747=\" },\n    @{ Old = \"42\u00d717=\"; New = \"14\u00d797=\" },\n    @{ Old = \"88\u00d784=\"; New = \"86\u00d762=\" },\n    @{ Old = \"55\u00d774=\"; New = \"82\u00d723=\" },\n    @{ Old = \"12\u00d793=\"; New = \"97\u00d768=\" },\n    @{ Old = \"72\u00d760=\"; New = \"90\u00d797=\" },\n    @{ Old = \"27\u00d728=\"; New = \"83\u00d781=\" },\n    @{ Old = \"61\u00d755=\"; New = \"59\u00d750=\" },\n    @{ Old = \"74\u00d736=\"; New = \"66\u00d786=\" },\n    @{ Old = \"36\u00d797=\"; New = \"97\u00d748=\" },\n    @{ Old = \"30\u00d711=\"; New = \"78\u00d744=\" },\n    @{ Old = \"44\u00d720=\"; New = \"93\u00d717=\" },\n    @{ Old = \"72\u00d755=\"; New = \"86\u00d712=\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n}\n"}
